$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 1
    3  = 4
    4  = 1
    5  = 1
    6  = 0
    7  = 3
    8  = 2
    9  = 2
    10 = 0
    11 = 0
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
